$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.64366028148607
$ws.Range("C2").Value = 17.36128077017902
$ws.Range("D2").Value = 17.88084631043447

$ws.Range("B3").Value = 1.776753170719777
$ws.Range("C3").Value = 2.020742317102819
$ws.Range("D3").Value = 2.519333872766607

$ws.Range("B4").Value = 0.3878329371642422
$ws.Range("C4").Value = 0.44855713502233
$ws.Range("D4").Value = 0.5557581910406075

$ws.Range("B5").Value = 79.11412320956629
$ws.Range("C5").Value = 79.87706986579832
$ws.Range("D5").Value = 80.8823178391079
